# Generate Report for Handoff
# Updates the localization-status report: refreshes file UUIDs / hashes /
# status / timestamps, and drops the now-unused "Latest Target File" (F)
# and "Latest Handback File" (G) columns' data from the per-language rows.

$wb = $excel.ActiveWorkbook

# ---- token substitutions (old -> new identifiers baked into file names / urls) ----
$uuidA_old = "5077eaab-51ab-4868-9300-0c7db760429e"
$uuidB_old = "863e2609-e7a8-4211-a0aa-5603a3d9c989"
$uuidA_new = "6f3021c9-7a13-43dc-9ed8-93935cb93275"
$uuidB_new = "ffffe047782c-e97e-4485-8c18-ac5f344f508d"

$hashA_old = "d36459d7f03bf69771da7f073578b52c6259f240"
$hashB_old = "b19cb24c14e837881260e26bd3365eddca6c3b0d"
$hash_new  = "ec89bef60032987ec7d9f825313fe9e7d75db401"

$status_old = "Handed back: in sync with en-US"
$status_new = "Ready for handoff"

function Convert-Token([string]$text) {
    $t = $text
    $t = $t.Replace($uuidA_old, $uuidA_new)
    $t = $t.Replace($uuidB_old, $uuidB_new)
    $t = $t.Replace($hashA_old, $hash_new)
    $t = $t.Replace($hashB_old, $hash_new)
    return $t
}

function Set-HyperlinkDisplayAndAddress($ws, [string]$addr, [string]$display, [string]$newAddress) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $display
            if ($newAddress) { $h.Address = $newAddress }
            return
        }
    }
}

function Remove-HyperlinkAt($ws, [string]$addr) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.Delete()
            return
        }
    }
}

# =====================================================================
# Sheet "Overview"
# =====================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $uuidA_new + ".md"
$wsOverview.Range("B2").Value = $status_new
$wsOverview.Range("C2").Value = $status_new
$wsOverview.Range("D2").Value = "2016-43-11 14:43:45"

$wsOverview.Range("A3").Value = $uuidB_new + ".md"
$wsOverview.Range("B3").Value = $status_new
$wsOverview.Range("C3").Value = $status_new
$wsOverview.Range("D3").Value = "2016-43-11 14:43:45"

$addrOvA2 = $null; $addrOvA3 = $null
foreach ($h in $wsOverview.Hyperlinks) {
    $a = $h.Range.Address()
    if ($a -eq "`$A`$2") { $addrOvA2 = $h.Address }
    elseif ($a -eq "`$A`$3") { $addrOvA3 = $h.Address }
}

Set-HyperlinkDisplayAndAddress $wsOverview "`$A`$2" ($uuidA_new + ".md") (Convert-Token $addrOvA2)
Set-HyperlinkDisplayAndAddress $wsOverview "`$A`$3" ($uuidB_new + ".md") (Convert-Token $addrOvA3)

# =====================================================================
# Helper that performs the shared per-language-sheet transform
# (sheet "zh-cn" and sheet "de-de" follow the identical row shape)
# =====================================================================
function Update-LanguageSheet($ws, [string]$handoffTime2, [string]$handoffTime3) {

    # --- capture existing hyperlink addresses (pre-edit) so we can token-substitute them ---
    $addrA2 = $null; $addrB2 = $null; $addrD2 = $null
    $addrA3 = $null; $addrB3 = $null; $addrD3 = $null
    foreach ($h in $ws.Hyperlinks) {
        $a = $h.Range.Address()
        if ($a -eq "`$A`$2") { $addrA2 = $h.Address }
        elseif ($a -eq "`$B`$2") { $addrB2 = $h.Address }
        elseif ($a -eq "`$D`$2") { $addrD2 = $h.Address }
        elseif ($a -eq "`$A`$3") { $addrA3 = $h.Address }
        elseif ($a -eq "`$B`$3") { $addrB3 = $h.Address }
        elseif ($a -eq "`$D`$3") { $addrD3 = $h.Address }
    }

    # --- remove the F/G (Latest Target File / Latest Handback File) cells + links ---
    Remove-HyperlinkAt $ws "`$F`$2"
    Remove-HyperlinkAt $ws "`$G`$2"
    Remove-HyperlinkAt $ws "`$F`$3"
    Remove-HyperlinkAt $ws "`$G`$3"
    $ws.Range("F2").Clear()
    $ws.Range("G2").Clear()
    $ws.Range("F3").Clear()
    $ws.Range("G3").Clear()

    # --- row 2 (uuid A file) ---
    $ws.Range("A2").Value = $uuidA_new + ".md"
    $ws.Range("B2").Value = ".md"
    $ws.Range("C2").Value = $status_new
    $ws.Range("D2").Value = (Convert-Token $ws.Range("D2").Value2)
    $ws.Range("E2").Value = $handoffTime2
    $ws.Range("H2").Value = "0001-01-01 00:00:00"
    $ws.Range("I2").Value = "Include"

    # --- row 3 (uuid B file) - D3 now references the SAME merged handoff file as row 2 ---
    $ws.Range("A3").Value = $uuidB_new + ".md"
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = $status_new
    $ws.Range("D3").Value = $ws.Range("D2").Value2
    $ws.Range("E3").Value = $handoffTime3
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("I3").Value = "Include"

    # --- refresh surviving hyperlinks (display text + target URL) ---
    Set-HyperlinkDisplayAndAddress $ws "`$A`$2" ($uuidA_new + ".md") (Convert-Token $addrA2)
    Set-HyperlinkDisplayAndAddress $ws "`$B`$2" ".md" (Convert-Token $addrB2)
    Set-HyperlinkDisplayAndAddress $ws "`$D`$2" ($ws.Range("D2").Value2) (Convert-Token $addrD2)

    Set-HyperlinkDisplayAndAddress $ws "`$A`$3" ($uuidB_new + ".md") (Convert-Token $addrA3)
    Set-HyperlinkDisplayAndAddress $ws "`$B`$3" ".md" (Convert-Token $addrB3)
    Set-HyperlinkDisplayAndAddress $ws "`$D`$3" ($ws.Range("D3").Value2) (Convert-Token $addrD3)
}

# =====================================================================
# Sheet "zh-cn"
# =====================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")
Update-LanguageSheet $wsZh "2016-03-11 14:43:42" "2016-03-11 14:43:42"

# =====================================================================
# Sheet "de-de"
# =====================================================================
$wsDe = $wb.Worksheets.Item("de-de")
Update-LanguageSheet $wsDe "2016-03-11 14:43:45" "2016-03-11 14:43:45"
